$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "Sprints" sheet - Sprint 1 now covers drafting the SRS (use case diagrams),
# Sprint 2 covers drafting the SRS system requirements. Sprint 1 is closed,
# Sprint 2 is active, and the sprint dates have moved to the new school year.
# ---------------------------------------------------------------------------
$sprints = $wb.Worksheets.Item("Sprints")

$sprints.Range("F4").Value = "Draft SRS, Use case diagrams"
$sprints.Range("G4").Value = 45412
$sprints.Range("I4").Value = "Closed"

$sprints.Range("F5").Value = "Draft SRS, System requirements for functional requirements"
$sprints.Range("I5").Value = "Active"

# ---------------------------------------------------------------------------
# "Product Backlog" sheet - the first story now targets the Library System
# (instead of the Automated Checkout System) and is Done; a second story for
# drafting the System Requirements has been added.
# ---------------------------------------------------------------------------
$backlog = $wb.Worksheets.Item("Product Backlog")

$backlog.Range("C5").Value = "Create System Requirements in SRS document"
$backlog.Range("E5").Value = "All necessary System Requirements are defined in the SRS document " + [char]10 + "SRS document is commited and pushed to Github"
$backlog.Range("D5").Value = "As a Software Developer, I need to define all the different System Requirements for the Library System"
$backlog.Range("F5").Value = 3
$backlog.Range("H5").Value = "Doing"
$backlog.Range("I5").Value = "Sprint 2"

$backlog.Range("D4").Value = "As a Software Developer, I need to define all the different Use Case Diagrams for the Library System"
$backlog.Range("H4").Value = "Done"

$wb.Save()
